# Actualización desde MV -datos-
# Appends two new daily rows (03-08-2021 and 04-08-2021) to the bottom of
# the "Derivados posiciones netas" table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $cellRef as a genuine text (shared-string) value.
# A plain "$ws.Range($cellRef).Value = $text" would let Excel's smart
# data-entry parser reinterpret strings such as "03-08-2021" as a date
# serial (since day=03 is also a valid month number), which would wrongly
# add a date number-format style to the sheet. Routing the text through a
# quoted formula and then pasting back as a value keeps it a literal string
# (matching every other date-label cell already in column A) without
# touching any cell's number format/style.
function Set-TextValue($cellRef, [string]$text) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

Set-TextValue "A148" "03-08-2021"
$ws.Range("B148").Value = 13080
$ws.Range("C148").Value = 21044
$ws.Range("D148").Value = -7963

Set-TextValue "A149" "04-08-2021"
$ws.Range("B149").Value = 13193
$ws.Range("C149").Value = 21223
$ws.Range("D149").Value = -8031
